$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 15).Value = "2022-07-25 20:57:52"
# Row 3
$ws.Cells.Item(3, 1).Value = "6075745012"
$ws.Cells.Item(3, 2).Value = "Avela Strumpfhose Madame Natural  11 - 12"
$ws.Cells.Item(3, 3).Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-madame-natural-11-12/p/6075745012"
$ws.Cells.Item(3, 8).Value = "5.95"
$ws.Cells.Item(3, 14).Value = "Avela Strumpfhose Madame Natural  11 - 12 5.95 Schweizer Franken"
$ws.Cells.Item(3, 15).Value = "2022-07-25 20:57:52"
# Row 4
$ws.Cells.Item(4, 1).Value = "6075749003"
$ws.Cells.Item(4, 2).Value = "Avela Strumpfhose Top Size Noir  11 - 12"
$ws.Cells.Item(4, 3).Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-top-size-noir-11-12/p/6075749003"
$ws.Cells.Item(4, 7).Value = "Avela"
$ws.Cells.Item(4, 8).Value = "1.95"
$ws.Cells.Item(4, 13).Value = "['haushalt-tier', 'bekleidung', 'socken-unterwaesche', 'struempfe']"
$ws.Cells.Item(4, 14).Value = "Avela Strumpfhose Top Size Noir  11 - 12 50% Aktion 1.95 Schweizer Franken statt 3.95 Schweizer Franken"
$ws.Cells.Item(4, 15).Value = "2022-07-25 20:57:52"
# Row 5
$ws.Cells.Item(5, 1).Value = "6031467017"
$ws.Cells.Item(5, 2).Value = "Naturaline Herren T-Shirt Kurzarm schwarz M"
$ws.Cells.Item(5, 3).Value = "/de/haushalt-tier/bekleidung/shirts-pullover/herren-shirt/naturaline-herren-t-shirt-kurzarm-schwarz-m/p/6031467017"
$ws.Cells.Item(5, 7).Value = "Coop"
$ws.Cells.Item(5, 8).Value = "24.95"
$ws.Cells.Item(5, 13).Value = "['haushalt-tier', 'bekleidung', 'shirts-pullover', 'herren-shirt']"
$ws.Cells.Item(5, 14).Value = "Naturaline Herren T-Shirt Kurzarm schwarz M 24.95 Schweizer Franken"
$ws.Cells.Item(5, 15).Value = "2022-07-25 20:57:52"
# Row 6
$ws.Cells.Item(6, 1).Value = "6075745013"
$ws.Cells.Item(6, 2).Value = "Avela Strumpfhose Madame Noir  8.5 - 9"
$ws.Cells.Item(6, 3).Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-madame-noir-85-9/p/6075745013"
$ws.Cells.Item(6, 5).Value = ""
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = "Avela"
$ws.Cells.Item(6, 8).Value = "5.95"
$ws.Cells.Item(6, 13).Value = "['haushalt-tier', 'bekleidung', 'socken-unterwaesche', 'struempfe']"
$ws.Cells.Item(6, 14).Value = "Avela Strumpfhose Madame Noir  8.5 - 9 5.95 Schweizer Franken"
$ws.Cells.Item(6, 15).Value = "2022-07-25 20:57:52"
# Row 7
$ws.Cells.Item(7, 1).Value = "3875554005"
$ws.Cells.Item(7, 2).Value = "Naturaline Damen Bustier Schwarz M"
$ws.Cells.Item(7, 3).Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-damen-bustier-schwarz-m/p/3875554005"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 5
$ws.Cells.Item(7, 8).Value = "14.95"
$ws.Cells.Item(7, 13).Value = "['haushalt-tier', 'bekleidung', 'socken-unterwaesche', 'unterwaesche']"
$ws.Cells.Item(7, 14).Value = "Naturaline Damen Bustier Schwarz M 14.95 Schweizer Franken"
$ws.Cells.Item(7, 15).Value = "2022-07-25 20:57:52"
# Row 8
$ws.Cells.Item(8, 1).Value = "6031467009"
$ws.Cells.Item(8, 2).Value = "Naturaline Herren T-Shirt Kurzarm weissXL"
$ws.Cells.Item(8, 3).Value = "/de/haushalt-tier/bekleidung/shirts-pullover/herren-shirt/naturaline-herren-t-shirt-kurzarm-weissxl/p/6031467009"
$ws.Cells.Item(8, 5).Value = ""
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 8).Value = "24.95"
$ws.Cells.Item(8, 13).Value = "['haushalt-tier', 'bekleidung', 'shirts-pullover', 'herren-shirt']"
$ws.Cells.Item(8, 14).Value = "Naturaline Herren T-Shirt Kurzarm weissXL - Online kein Bestand 24.95 Schweizer Franken"
$ws.Cells.Item(8, 15).Value = "2022-07-25 20:57:52"
# Row 9
$ws.Cells.Item(9, 1).Value = "3875554009"
$ws.Cells.Item(9, 2).Value = "Naturaline Damen Bustier Weiss L"
$ws.Cells.Item(9, 3).Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-damen-bustier-weiss-l/p/3875554009"
$ws.Cells.Item(9, 6).Value = 5
$ws.Cells.Item(9, 8).Value = "14.95"
$ws.Cells.Item(9, 14).Value = "Naturaline Damen Bustier Weiss L 14.95 Schweizer Franken"
$ws.Cells.Item(9, 15).Value = "2022-07-25 20:57:52"
# Row 10
$ws.Cells.Item(10, 1).Value = "3305779007"
$ws.Cells.Item(10, 2).Value = "Naturaline Damen Panty S weiss"
$ws.Cells.Item(10, 3).Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-damen-panty-s-weiss/p/3305779007"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 4
$ws.Cells.Item(10, 14).Value = "Naturaline Damen Panty S weiss 9.95 Schweizer Franken"
$ws.Cells.Item(10, 15).Value = "2022-07-25 20:57:52"
# Row 11
$ws.Cells.Item(11, 1).Value = "3305289015"
$ws.Cells.Item(11, 2).Value = "Naturaline Herren Slip weiss L"
$ws.Cells.Item(11, 3).Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-herren-slip-weiss-l/p/3305289015"
$ws.Cells.Item(11, 13).Value = "['haushalt-tier', 'bekleidung', 'socken-unterwaesche', 'unterwaesche']"
$ws.Cells.Item(11, 14).Value = "Naturaline Herren Slip weiss L 9.95 Schweizer Franken"
$ws.Cells.Item(11, 15).Value = "2022-07-25 20:57:52"
# Row 12
$ws.Cells.Item(12, 1).Value = "6365980012"
$ws.Cells.Item(12, 2).Value = "Naturaline Herren Socken Glatt Duo Weiss 43 - 45"
$ws.Cells.Item(12, 3).Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/naturaline-herren-socken-glatt-duo-weiss-43-45/p/6365980012"
$ws.Cells.Item(12, 8).Value = "9.95"
$ws.Cells.Item(12, 13).Value = "['haushalt-tier', 'bekleidung', 'socken-unterwaesche', 'socken']"
$ws.Cells.Item(12, 14).Value = "Naturaline Herren Socken Glatt Duo Weiss 43 - 45 9.95 Schweizer Franken"
$ws.Cells.Item(12, 15).Value = "2022-07-25 20:57:52"
# Row 13
$ws.Cells.Item(13, 1).Value = "6031467006"
$ws.Cells.Item(13, 2).Value = "Naturaline Herren T-Shirt Kurzarm weiss S"
$ws.Cells.Item(13, 3).Value = "/de/haushalt-tier/bekleidung/shirts-pullover/herren-shirt/naturaline-herren-t-shirt-kurzarm-weiss-s/p/6031467006"
$ws.Cells.Item(13, 4).Value = ""
$ws.Cells.Item(13, 7).Value = "Coop"
$ws.Cells.Item(13, 8).Value = "24.95"
$ws.Cells.Item(13, 9).Value = ""
$ws.Cells.Item(13, 10).Value = ""
$ws.Cells.Item(13, 11).Value = ""
$ws.Cells.Item(13, 12).Value = ""
$ws.Cells.Item(13, 13).Value = "['haushalt-tier', 'bekleidung', 'shirts-pullover', 'herren-shirt']"
$ws.Cells.Item(13, 14).Value = "Naturaline Herren T-Shirt Kurzarm weiss S 24.95 Schweizer Franken"
$ws.Cells.Item(13, 15).Value = "2022-07-25 20:57:52"
# Row 14
$ws.Cells.Item(14, 1).Value = "6077158006"
$ws.Cells.Item(14, 2).Value = "Avela Söckchen Ideal Noir One Size"
$ws.Cells.Item(14, 3).Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/avela-soeckchen-ideal-noir-one-size/p/6077158006"
$ws.Cells.Item(14, 4).Value = "2ST"
$ws.Cells.Item(14, 8).Value = "4.95"
$ws.Cells.Item(14, 9).Value = "2.48/1ST"
$ws.Cells.Item(14, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(14, 11).Value = "2.48"
$ws.Cells.Item(14, 12).Value = "1ST"
$ws.Cells.Item(14, 13).Value = "['haushalt-tier', 'bekleidung', 'socken-unterwaesche', 'socken']"
$ws.Cells.Item(14, 14).Value = "Avela Söckchen Ideal Noir One Size 4.95 Schweizer Franken"
$ws.Cells.Item(14, 15).Value = "2022-07-25 20:57:52"
# Row 15
$ws.Cells.Item(15, 1).Value = "6075745006"
$ws.Cells.Item(15, 2).Value = "Avela Strumpfhose Madame Hasel  9.5"
$ws.Cells.Item(15, 3).Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-madame-hasel-95/p/6075745006"
$ws.Cells.Item(15, 8).Value = "5.95"
$ws.Cells.Item(15, 14).Value = "Avela Strumpfhose Madame Hasel  9.5 5.95 Schweizer Franken"
$ws.Cells.Item(15, 15).Value = "2022-07-25 20:57:52"
# Row 16
$ws.Cells.Item(16, 15).Value = "2022-07-25 20:57:52"
# Row 17
$ws.Cells.Item(17, 15).Value = "2022-07-25 20:57:52"
# Row 18
$ws.Cells.Item(18, 15).Value = "2022-07-25 20:57:52"
# Row 19
$ws.Cells.Item(19, 15).Value = "2022-07-25 20:57:52"
# Row 20
$ws.Cells.Item(20, 15).Value = "2022-07-25 20:57:52"
# Row 21
$ws.Cells.Item(21, 15).Value = "2022-07-25 20:57:52"
# Row 22
$ws.Cells.Item(22, 15).Value = "2022-07-25 20:57:52"
# Row 23
$ws.Cells.Item(23, 15).Value = "2022-07-25 20:57:52"
# Row 24
$ws.Cells.Item(24, 15).Value = "2022-07-25 20:57:52"
# Row 25
$ws.Cells.Item(25, 15).Value = "2022-07-25 20:57:52"
# Row 26
$ws.Cells.Item(26, 15).Value = "2022-07-25 20:57:52"
# Row 27
$ws.Cells.Item(27, 15).Value = "2022-07-25 20:57:52"
# Row 28
$ws.Cells.Item(28, 15).Value = "2022-07-25 20:57:52"
# Row 29
$ws.Cells.Item(29, 15).Value = "2022-07-25 20:57:52"
# Row 30
$ws.Cells.Item(30, 15).Value = "2022-07-25 20:57:52"
# Row 31
$ws.Cells.Item(31, 15).Value = "2022-07-25 20:57:52"
# Row 32
$ws.Cells.Item(32, 15).Value = "2022-07-25 20:57:52"
# Row 33
$ws.Cells.Item(33, 15).Value = "2022-07-25 20:57:52"
# Row 34
$ws.Cells.Item(34, 15).Value = "2022-07-25 20:57:52"
# Row 35
$ws.Cells.Item(35, 15).Value = "2022-07-25 20:57:52"
# Row 36
$ws.Cells.Item(36, 15).Value = "2022-07-25 20:57:52"
# Row 37
$ws.Cells.Item(37, 15).Value = "2022-07-25 20:57:52"
# Row 38
$ws.Cells.Item(38, 15).Value = "2022-07-25 20:57:52"
# Row 39
$ws.Cells.Item(39, 15).Value = "2022-07-25 20:57:52"
# Row 40
$ws.Cells.Item(40, 15).Value = "2022-07-25 20:57:52"
# Row 41
$ws.Cells.Item(41, 15).Value = "2022-07-25 20:57:52"
# Row 42
$ws.Cells.Item(42, 15).Value = "2022-07-25 20:57:52"
# Row 43
$ws.Cells.Item(43, 15).Value = "2022-07-25 20:57:52"
# Row 44
$ws.Cells.Item(44, 15).Value = "2022-07-25 20:57:52"
# Row 45
$ws.Cells.Item(45, 15).Value = "2022-07-25 20:57:52"
# Row 46
$ws.Cells.Item(46, 15).Value = "2022-07-25 20:57:52"
# Row 47
$ws.Cells.Item(47, 15).Value = "2022-07-25 20:57:52"
# Row 48
$ws.Cells.Item(48, 15).Value = "2022-07-25 20:57:52"
# Row 49
$ws.Cells.Item(49, 15).Value = "2022-07-25 20:57:52"
# Row 50
$ws.Cells.Item(50, 15).Value = "2022-07-25 20:57:52"
# Row 51
$ws.Cells.Item(51, 15).Value = "2022-07-25 20:57:52"
# Row 52
$ws.Cells.Item(52, 15).Value = "2022-07-25 20:57:52"
# Row 53
$ws.Cells.Item(53, 15).Value = "2022-07-25 20:57:52"
# Row 54
$ws.Cells.Item(54, 15).Value = "2022-07-25 20:57:52"
# Row 55
$ws.Cells.Item(55, 15).Value = "2022-07-25 20:57:52"
# Row 56
$ws.Cells.Item(56, 15).Value = "2022-07-25 20:57:52"
# Row 57
$ws.Cells.Item(57, 15).Value = "2022-07-25 20:57:52"
# Row 58
$ws.Cells.Item(58, 15).Value = "2022-07-25 20:57:52"
# Row 59
$ws.Cells.Item(59, 15).Value = "2022-07-25 20:57:52"
# Row 60
$ws.Cells.Item(60, 15).Value = "2022-07-25 20:57:52"
# Row 61
$ws.Cells.Item(61, 15).Value = "2022-07-25 20:57:52"
# Row 62
$ws.Cells.Item(62, 15).Value = "2022-07-25 20:57:52"
# Row 63
$ws.Cells.Item(63, 15).Value = "2022-07-25 20:57:52"
# Row 64
$ws.Cells.Item(64, 15).Value = "2022-07-25 20:57:52"
# Row 65
$ws.Cells.Item(65, 15).Value = "2022-07-25 20:57:52"
# Row 66
$ws.Cells.Item(66, 15).Value = "2022-07-25 20:57:52"
# Row 67
$ws.Cells.Item(67, 15).Value = "2022-07-25 20:57:52"
# Row 68
$ws.Cells.Item(68, 15).Value = "2022-07-25 20:57:52"
# Row 69
$ws.Cells.Item(69, 15).Value = "2022-07-25 20:57:52"
# Row 70
$ws.Cells.Item(70, 15).Value = "2022-07-25 20:57:52"
# Row 71
$ws.Cells.Item(71, 15).Value = "2022-07-25 20:57:52"
# Row 72
$ws.Cells.Item(72, 15).Value = "2022-07-25 20:57:52"
# Row 73
$ws.Cells.Item(73, 15).Value = "2022-07-25 20:57:52"
